# Remove column from alcohol data:
# The sheet had an extra column (M) that duplicated/derived data; the
# correct data (previously in column N) shifts left to become the new
# column M, and the old column M is discarded entirely.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Delete the whole column M on the data sheet; this shifts column N (and
# beyond) one position to the left, so the former N values become M.
$ws1.Columns.Item(13).Delete()

# Reselect the data sheet, zoom in and move the selection to the new
# right-most column (M1).
$ws1.Select()
$excel.ActiveWindow.Zoom = 110
$null = $ws1.Range("M1").Select()

# The other (empty) sheets were also re-zoomed when the workbook was saved.
$ws2.Select()
$excel.ActiveWindow.Zoom = 110

$ws3.Select()
$excel.ActiveWindow.Zoom = 110

# Leave the original data sheet as the active / selected sheet.
$ws1.Select()
